# Refresh the scraped cryptocurrency Price (D) / Volume(1h) (E) figures, and
# fix the Hedera / TrustWalletToken row ordering, per the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.366.66'
$ws.Range('E2').Value = '  -3.90%  '
$ws.Range('D3').Value = '1.862.26'
$ws.Range('E3').Value = '  -4.72%  '
$ws.Range('E4').Value = '  -1.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4541'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3875'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.20'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -10.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07932'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.67%  '
$ws.Range('E11').Value = '  -3.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.43'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.54%  '
$ws.Range('D13').Value = '1.837.37'
$ws.Range('E13').Value = '  -5.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.910'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.08%  '
$ws.Range('E15').Value = '  -5.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9998'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001036'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '85.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06512'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.536'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.20%  '
$ws.Range('D23').Value = '27.357.04'
$ws.Range('E23').Value = '  -3.96%  '
$ws.Range('E24').Value = '  -4.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.275'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('D26').Value = '2.076.19'
$ws.Range('E26').Value = '  -4.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '152.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.24%  '
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.067'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.503'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '120.89'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.495'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09337'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9360'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.610'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('E36').Value = '  -6.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02241'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.12%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06004'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.09%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.222'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.280'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.39%  '
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5919'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1890'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.16'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.32%  '
$ws.Range('E45').Value = '  -3.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5633'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.930'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.373'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06792'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.19'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.58%  '
